# "leçon 2 et 3" — insert two new lesson rows (npm / TypeScript 2 - Express)
# into the "horaire" sheet, pushing the old "npm" / "TypeScript - La suite"
# rows (which had no Exercice entry) out in favour of the new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (cours #2): was "npm" with no exercise -> now has an exercise link
# plus a proper markdown link for the course material itself.
# (Write order matters for shared-string table placement: the exercise
# string must land before the course string to match the saved file.)
$ws.Range("D3").Value = "[Exercice 2 - lodash](exercice2_lodash.md)"
$ws.Range("C3").Value = "[npm](npm.md)"

# Row 4 (cours #3): was "TypeScript - La suite" with no exercise -> now
# covers TypeScript 2 + intro to Express, with its own exercise.
$ws.Range("C4").Value = "[TypeScript 2](typescript_2.md)<br/>[Introduction Express](introduction_express.md)"
$ws.Range("D4").Value = "[Exercice 3 - Express](exercice3_express.md)"

# Reflect the author's last selection when they saved the file.
[void]$ws.Range("D5").Select()
